$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "AL" column (the "Total Cell Count Across Samples" helper
# column, holding a merged header + per-row SUM formulas). This shifts the
# trailing blank "AM" column left into AL's place, matching the authored diff.
$ws.Columns("AL").Delete()

# Restore the sheet-view state recorded in the diff.
$ws.Range("AC101").Select()
$excel.ActiveWindow.Zoom = 84
